$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "43÷8=5, 3"
$tbl.Cell(1, 2).Range.Text = "84÷6=14, 0"
$tbl.Cell(1, 3).Range.Text = "62÷8=7, 6"
$tbl.Cell(1, 4).Range.Text = "18÷6=3, 0"
$tbl.Cell(1, 5).Range.Text = "30÷5=6, 0"

$tbl.Cell(5, 1).Range.Text = "11÷3=3, 2"
$tbl.Cell(5, 2).Range.Text = "43÷2=21, 1"
$tbl.Cell(5, 3).Range.Text = "61÷3=20, 1"
$tbl.Cell(5, 4).Range.Text = "44÷6=7, 2"
$tbl.Cell(5, 5).Range.Text = "10÷6=1, 4"

$tbl.Cell(9, 1).Range.Text = "92÷3=30, 2"
$tbl.Cell(9, 2).Range.Text = "61÷3=20, 1"
$tbl.Cell(9, 3).Range.Text = "27÷4=6, 3"
$tbl.Cell(9, 4).Range.Text = "99÷5=19, 4"
$tbl.Cell(9, 5).Range.Text = "45÷4=11, 1"

$tbl.Cell(13, 1).Range.Text = "85÷2=42, 1"
$tbl.Cell(13, 2).Range.Text = "46÷8=5, 6"
$tbl.Cell(13, 3).Range.Text = "60÷4=15, 0"
$tbl.Cell(13, 4).Range.Text = "46÷3=15, 1"
$tbl.Cell(13, 5).Range.Text = "62÷7=8, 6"

$tbl.Cell(17, 1).Range.Text = "77÷8=9, 5"
$tbl.Cell(17, 2).Range.Text = "40÷2=20, 0"
$tbl.Cell(17, 3).Range.Text = "54÷6=9, 0"
$tbl.Cell(17, 4).Range.Text = "91÷2=45, 1"
$tbl.Cell(17, 5).Range.Text = "81÷4=20, 1"
